$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly price-record data between row 4 and row 5
# (columns: Fecha, Calidad, Volumen, Precio minimo, Precio maximo,
#  Precio promedio ponderado, Unidad de comercializacion, Precio $/Kg, Kg o Unidades)

$row4_D = $ws.Range("D4").Value2
$row4_I = $ws.Range("I4").Value2
$row4_J = $ws.Range("J4").Value2
$row4_K = $ws.Range("K4").Value2
$row4_L = $ws.Range("L4").Value2
$row4_M = $ws.Range("M4").Value2
$row4_N = $ws.Range("N4").Value2
$row4_P = $ws.Range("P4").Value2
$row4_Q = $ws.Range("Q4").Value2

$row5_D = $ws.Range("D5").Value2
$row5_I = $ws.Range("I5").Value2
$row5_J = $ws.Range("J5").Value2
$row5_K = $ws.Range("K5").Value2
$row5_L = $ws.Range("L5").Value2
$row5_M = $ws.Range("M5").Value2
$row5_N = $ws.Range("N5").Value2
$row5_P = $ws.Range("P5").Value2
$row5_Q = $ws.Range("Q5").Value2

$ws.Range("D4").Value = $row5_D
$ws.Range("I4").Value = $row5_I
$ws.Range("J4").Value = $row5_J
$ws.Range("K4").Value = $row5_K
$ws.Range("L4").Value = $row5_L
$ws.Range("M4").Value = $row5_M
$ws.Range("N4").Value = $row5_N
$ws.Range("P4").Value = $row5_P
$ws.Range("Q4").Value = $row5_Q

$ws.Range("D5").Value = $row4_D
$ws.Range("I5").Value = $row4_I
$ws.Range("J5").Value = $row4_J
$ws.Range("K5").Value = $row4_K
$ws.Range("L5").Value = $row4_L
$ws.Range("M5").Value = $row4_M
$ws.Range("N5").Value = $row4_N
$ws.Range("P5").Value = $row4_P
$ws.Range("Q5").Value = $row4_Q
